# Natmi following Dr Hou advice
# Expand the Wnt2 -> Fzd3 sender/target-cluster matrix from 3 rows (one sender x
# one target each) to the full 2 senders (ECs, FAPs) x 5 targets (ECs, FAPs, M2,
# Neutro, sCs) cross-product, with recomputed ligand/receptor/edge statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs (ligand Wnt2 / receptor Fzd3)
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Wnt2"
$ws.Cells.Item(2,3).Value = "Fzd3"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 1.0
$ws.Cells.Item(2,6).Value = 0.3333333333333333
$ws.Cells.Item(2,7).Value = 0.01070233333333333
$ws.Cells.Item(2,8).Value = 0.032107
$ws.Cells.Item(2,9).Value = 0.006017198313602724
$ws.Cells.Item(2,10).Value = 0.006017198313602724
$ws.Cells.Item(2,11).Value = 3.0
$ws.Cells.Item(2,12).Value = 1.0
$ws.Cells.Item(2,13).Value = 0.174733
$ws.Cells.Item(2,14).Value = 0.524199
$ws.Cells.Item(2,15).Value = 0.05882867859784573
$ws.Cells.Item(2,16).Value = 0.05882867859784573
$ws.Cells.Item(2,17).Value = 0.001870050810333333
$ws.Cells.Item(2,18).Value = 0.016830457293
$ws.Cells.Item(2,19).Value = 0.000353983825650434
$ws.Cells.Item(2,20).Value = 0.000353983825650434

# Row 3: ECs -> FAPs (ligand Wnt2 / receptor Fzd3)
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Wnt2"
$ws.Cells.Item(3,3).Value = "Fzd3"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 1.0
$ws.Cells.Item(3,6).Value = 0.3333333333333333
$ws.Cells.Item(3,7).Value = 0.01070233333333333
$ws.Cells.Item(3,8).Value = 0.032107
$ws.Cells.Item(3,9).Value = 0.006017198313602724
$ws.Cells.Item(3,10).Value = 0.006017198313602724
$ws.Cells.Item(3,11).Value = 3.0
$ws.Cells.Item(3,12).Value = 1.0
$ws.Cells.Item(3,13).Value = 0.4828523333333334
$ws.Cells.Item(3,14).Value = 1.448557
$ws.Cells.Item(3,15).Value = 0.1625655412995058
$ws.Cells.Item(3,16).Value = 0.1625655412995058
$ws.Cells.Item(3,17).Value = 0.005167646622111111
$ws.Cells.Item(3,18).Value = 0.046508819599
$ws.Cells.Item(3,19).Value = 0.0009781891009573
$ws.Cells.Item(3,20).Value = 0.0009781891009573002

# Row 4: ECs -> M2 (ligand Wnt2 / receptor Fzd3)
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Wnt2"
$ws.Cells.Item(4,3).Value = "Fzd3"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 1.0
$ws.Cells.Item(4,6).Value = 0.3333333333333333
$ws.Cells.Item(4,7).Value = 0.01070233333333333
$ws.Cells.Item(4,8).Value = 0.032107
$ws.Cells.Item(4,9).Value = 0.006017198313602724
$ws.Cells.Item(4,10).Value = 0.006017198313602724
$ws.Cells.Item(4,11).Value = 1.0
$ws.Cells.Item(4,12).Value = 0.3333333333333333
$ws.Cells.Item(4,13).Value = 0.006173333333333333
$ws.Cells.Item(4,14).Value = 0.01852
$ws.Cells.Item(4,15).Value = 0.002078422750963094
$ws.Cells.Item(4,16).Value = 0.002078422750963094
$ws.Cells.Item(4,17).Value = 0.00006606907111111111
$ws.Cells.Item(4,18).Value = 0.0005946216399999999
$ws.Cells.Item(4,19).Value = 0.00001250628187204866
$ws.Cells.Item(4,20).Value = 0.00001250628187204866

# Row 5: ECs -> Neutro (ligand Wnt2 / receptor Fzd3)
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Wnt2"
$ws.Cells.Item(5,3).Value = "Fzd3"
$ws.Cells.Item(5,4).Value = "Neutro"
$ws.Cells.Item(5,5).Value = 1.0
$ws.Cells.Item(5,6).Value = 0.3333333333333333
$ws.Cells.Item(5,7).Value = 0.01070233333333333
$ws.Cells.Item(5,8).Value = 0.032107
$ws.Cells.Item(5,9).Value = 0.006017198313602724
$ws.Cells.Item(5,10).Value = 0.006017198313602724
$ws.Cells.Item(5,11).Value = 2.0
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.02610233333333334
$ws.Cells.Item(5,14).Value = 0.078307
$ws.Cells.Item(5,15).Value = 0.00878806967384811
$ws.Cells.Item(5,16).Value = 0.008788069673848112
$ws.Cells.Item(5,17).Value = 0.0002793558721111111
$ws.Cells.Item(5,18).Value = 0.002514202849
$ws.Cells.Item(5,19).Value = 0.00005287955802130209
$ws.Cells.Item(5,20).Value = 0.0000528795580213021

# Row 6: ECs -> sCs (ligand Wnt2 / receptor Fzd3)
$ws.Cells.Item(6,1).Value = "ECs"
$ws.Cells.Item(6,2).Value = "Wnt2"
$ws.Cells.Item(6,3).Value = "Fzd3"
$ws.Cells.Item(6,4).Value = "sCs"
$ws.Cells.Item(6,5).Value = 1.0
$ws.Cells.Item(6,6).Value = 0.3333333333333333
$ws.Cells.Item(6,7).Value = 0.01070233333333333
$ws.Cells.Item(6,8).Value = 0.032107
$ws.Cells.Item(6,9).Value = 0.006017198313602724
$ws.Cells.Item(6,10).Value = 0.006017198313602724
$ws.Cells.Item(6,11).Value = 3.0
$ws.Cells.Item(6,12).Value = 1.0
$ws.Cells.Item(6,13).Value = 2.28034
$ws.Cells.Item(6,14).Value = 6.84102
$ws.Cells.Item(6,15).Value = 0.7677392876778373
$ws.Cells.Item(6,16).Value = 0.7677392876778373
$ws.Cells.Item(6,17).Value = 0.02440495879333333
$ws.Cells.Item(6,18).Value = 0.21964462914
$ws.Cells.Item(6,19).Value = 0.004619639547101639
$ws.Cells.Item(6,20).Value = 0.004619639547101639

# Row 7: FAPs -> ECs (ligand Wnt2 / receptor Fzd3)
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Wnt2"
$ws.Cells.Item(7,3).Value = "Fzd3"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 3.0
$ws.Cells.Item(7,6).Value = 1.0
$ws.Cells.Item(7,7).Value = 1.767921666666667
$ws.Cells.Item(7,8).Value = 5.303765
$ws.Cells.Item(7,9).Value = 0.9939828016863973
$ws.Cells.Item(7,10).Value = 0.9939828016863973
$ws.Cells.Item(7,11).Value = 3.0
$ws.Cells.Item(7,12).Value = 1.0
$ws.Cells.Item(7,13).Value = 0.174733
$ws.Cells.Item(7,14).Value = 0.524199
$ws.Cells.Item(7,15).Value = 0.05882867859784573
$ws.Cells.Item(7,16).Value = 0.05882867859784573
$ws.Cells.Item(7,17).Value = 0.3089142565816667
$ws.Cells.Item(7,18).Value = 2.780228309235
$ws.Cells.Item(7,19).Value = 0.05847469477219529
$ws.Cells.Item(7,20).Value = 0.0584746947721953

# Row 8: FAPs -> FAPs (ligand Wnt2 / receptor Fzd3)
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Wnt2"
$ws.Cells.Item(8,3).Value = "Fzd3"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 3.0
$ws.Cells.Item(8,6).Value = 1.0
$ws.Cells.Item(8,7).Value = 1.767921666666667
$ws.Cells.Item(8,8).Value = 5.303765
$ws.Cells.Item(8,9).Value = 0.9939828016863973
$ws.Cells.Item(8,10).Value = 0.9939828016863973
$ws.Cells.Item(8,11).Value = 3.0
$ws.Cells.Item(8,12).Value = 1.0
$ws.Cells.Item(8,13).Value = 0.4828523333333334
$ws.Cells.Item(8,14).Value = 1.448557
$ws.Cells.Item(8,15).Value = 0.1625655412995058
$ws.Cells.Item(8,16).Value = 0.1625655412995058
$ws.Cells.Item(8,17).Value = 0.8536451019005558
$ws.Cells.Item(8,18).Value = 7.682805917105001
$ws.Cells.Item(8,19).Value = 0.1615873521985485
$ws.Cells.Item(8,20).Value = 0.1615873521985485

# Row 9: FAPs -> M2 (ligand Wnt2 / receptor Fzd3)
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Wnt2"
$ws.Cells.Item(9,3).Value = "Fzd3"
$ws.Cells.Item(9,4).Value = "M2"
$ws.Cells.Item(9,5).Value = 3.0
$ws.Cells.Item(9,6).Value = 1.0
$ws.Cells.Item(9,7).Value = 1.767921666666667
$ws.Cells.Item(9,8).Value = 5.303765
$ws.Cells.Item(9,9).Value = 0.9939828016863973
$ws.Cells.Item(9,10).Value = 0.9939828016863973
$ws.Cells.Item(9,11).Value = 1.0
$ws.Cells.Item(9,12).Value = 0.3333333333333333
$ws.Cells.Item(9,13).Value = 0.006173333333333333
$ws.Cells.Item(9,14).Value = 0.01852
$ws.Cells.Item(9,15).Value = 0.002078422750963094
$ws.Cells.Item(9,16).Value = 0.002078422750963094
$ws.Cells.Item(9,17).Value = 0.01091396975555556
$ws.Cells.Item(9,18).Value = 0.0982257278
$ws.Cells.Item(9,19).Value = 0.002065916469091045
$ws.Cells.Item(9,20).Value = 0.002065916469091046

# Row 10: FAPs -> Neutro (ligand Wnt2 / receptor Fzd3)
$ws.Cells.Item(10,1).Value = "FAPs"
$ws.Cells.Item(10,2).Value = "Wnt2"
$ws.Cells.Item(10,3).Value = "Fzd3"
$ws.Cells.Item(10,4).Value = "Neutro"
$ws.Cells.Item(10,5).Value = 3.0
$ws.Cells.Item(10,6).Value = 1.0
$ws.Cells.Item(10,7).Value = 1.767921666666667
$ws.Cells.Item(10,8).Value = 5.303765
$ws.Cells.Item(10,9).Value = 0.9939828016863973
$ws.Cells.Item(10,10).Value = 0.9939828016863973
$ws.Cells.Item(10,11).Value = 2.0
$ws.Cells.Item(10,12).Value = 0.6666666666666666
$ws.Cells.Item(10,13).Value = 0.02610233333333334
$ws.Cells.Item(10,14).Value = 0.078307
$ws.Cells.Item(10,15).Value = 0.00878806967384811
$ws.Cells.Item(10,16).Value = 0.008788069673848112
$ws.Cells.Item(10,17).Value = 0.04614688065055556
$ws.Cells.Item(10,18).Value = 0.4153219258550001
$ws.Cells.Item(10,19).Value = 0.008735190115826809
$ws.Cells.Item(10,20).Value = 0.00873519011582681

# Row 11: FAPs -> sCs (ligand Wnt2 / receptor Fzd3)
$ws.Cells.Item(11,1).Value = "FAPs"
$ws.Cells.Item(11,2).Value = "Wnt2"
$ws.Cells.Item(11,3).Value = "Fzd3"
$ws.Cells.Item(11,4).Value = "sCs"
$ws.Cells.Item(11,5).Value = 3.0
$ws.Cells.Item(11,6).Value = 1.0
$ws.Cells.Item(11,7).Value = 1.767921666666667
$ws.Cells.Item(11,8).Value = 5.303765
$ws.Cells.Item(11,9).Value = 0.9939828016863973
$ws.Cells.Item(11,10).Value = 0.9939828016863973
$ws.Cells.Item(11,11).Value = 3.0
$ws.Cells.Item(11,12).Value = 1.0
$ws.Cells.Item(11,13).Value = 2.28034
$ws.Cells.Item(11,14).Value = 6.84102
$ws.Cells.Item(11,15).Value = 0.7677392876778373
$ws.Cells.Item(11,16).Value = 0.7677392876778373
$ws.Cells.Item(11,17).Value = 4.031462493366668
$ws.Cells.Item(11,18).Value = 36.2831624403
$ws.Cells.Item(11,19).Value = 0.7631196481307356
$ws.Cells.Item(11,20).Value = 0.7631196481307356
